# Rename the three header/footer logo pictures' display names:
#   footer1.xml (Footers "first")   : image1.png -> image2.png
#   footer2.xml (Footers "default") : image1.png -> image2.png
#   header1.xml (Headers "first")   : image2.jpg -> image1.jpg
#
# The pictures live in the section's headers/footers (not the main body),
# so they are reached through Section.Headers / Section.Footers rather than
# ActiveDocument.InlineShapes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoPicture($headerFooter, $newName) {
    # Grab the picture, then re-fetch it through its own Range before
    # writing .Name - setting the property directly off the
    # HeaderFooter.Range collection can fail to resolve on some stories.
    $shape = $headerFooter.Range.InlineShapes.Item(1)
    $shapeRange = $shape.Range
    $shape = $shapeRange.InlineShapes.Item(1)
    $shape.Name = $newName
}

# footer1.xml == the "first page" footer == Footers.Item(2)
Rename-LogoPicture $sec.Footers.Item(2) "image2.png"

# footer2.xml == the "default" footer == Footers.Item(1)
Rename-LogoPicture $sec.Footers.Item(1) "image2.png"

# header1.xml == the "first page" header == Headers.Item(2)
Rename-LogoPicture $sec.Headers.Item(2) "image1.jpg"
